$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current header row (row 1). This pushes the
# existing header row ("Lg.", "Threading", ...) down to row 2, and all the
# data rows shift down by one as well.
$ws.Rows("1:1").Insert()

# Copy the header formatting (bold font, thin borders, centered/top aligned)
# from the old header row - now at row 2 - onto the brand-new row 1.
$ws.Range("A2:K2").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)   # xlPasteFormats

# The old header row (now row 2) should go back to plain/unstyled formatting.
$ws.Rows("2:2").ClearFormats()

# Populate the new row 1 with a simple numeric index (0 - 10).
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10

# On the (now second) header row, a few cells lose their text content.
$ws.Range("H2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
